$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tree table is being collapsed down to just the header row plus a
# single data row (columns C:D and rows 3:9 go away entirely).
$ws.Range("A3:D9").Delete(-4162)   # xlShiftUp -- drop the old data rows
$ws.Range("C1:D2").Delete(-4159)   # xlShiftToLeft -- drop the OPK2/OPK3 columns

# Refresh the remaining header labels (A1/B1 keep their existing bold/
# bordered header style -- only the text changes).
$ws.Range("A1").Value = "Drzewo"
$ws.Range("B1").Value = "OPK1"

# Row 2 becomes the sole data row. "73" must stay textual (not numeric),
# so the cell is pre-formatted as Text before the value is typed in, then
# reset back to the workbook's Normal style once the string is locked in.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "73"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "MO"
